$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Cells.Item(7, 1).Value = '67jtezzI'
$ws.Cells.Item(7, 2).Value = '24/10/2024'
$ws.Cells.Item(7, 3).Value = '14:00'
$ws.Cells.Item(7, 4).Value = 'SPAIN - LALIGA2'
$ws.Cells.Item(7, 5).Value = 'Almeria'
$ws.Cells.Item(7, 6).Value = 'Albacete'
$ws.Cells.Item(7, 7).Value = 1.62
$ws.Cells.Item(7, 8).Value = 3.75
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 2.2
$ws.Cells.Item(7, 11).Value = 2.5
$ws.Cells.Item(7, 12).Value = 4.5
$ws.Cells.Item(7, 13).Value = 1.03
$ws.Cells.Item(7, 14).Value = 17
$ws.Cells.Item(7, 15).Value = 1.17
$ws.Cells.Item(7, 16).Value = 5
$ws.Cells.Item(7, 17).Value = 1.53
$ws.Cells.Item(7, 18).Value = 2.4
$ws.Cells.Item(7, 19).Value = 1.29
$ws.Cells.Item(7, 20).Value = 3.5
$ws.Cells.Item(7, 21).Value = 1.57
$ws.Cells.Item(7, 22).Value = 2.25
$ws.Cells.Item(7, 23).Value = 10
$ws.Cells.Item(7, 24).Value = 9.5
$ws.Cells.Item(7, 25).Value = 8.5
$ws.Cells.Item(7, 26).Value = 13
$ws.Cells.Item(7, 27).Value = 12
$ws.Cells.Item(7, 28).Value = 21
$ws.Cells.Item(7, 29).Value = 15
$ws.Cells.Item(7, 30).Value = 7.5
$ws.Cells.Item(7, 31).Value = 12
$ws.Cells.Item(7, 32).Value = 41
$ws.Cells.Item(7, 33).Value = 126
$ws.Cells.Item(7, 34).Value = 19
$ws.Cells.Item(7, 35).Value = 29
$ws.Cells.Item(7, 36).Value = 17
$ws.Cells.Item(7, 37).Value = 51
$ws.Cells.Item(7, 38).Value = 34
$ws.Cells.Item(7, 39).Value = 34
$ws.Cells.Item(7, 40).Value = 4
$ws.Cells.Item(7, 41).Value = 8.5
$ws.Cells.Item(7, 42).Value = 15
$ws.Cells.Item(7, 43).Value = 23
$ws.Cells.Item(7, 44).Value = 41
$ws.Cells.Item(7, 45).Value = 81
$ws.Cells.Item(7, 46).Value = 3.5
$ws.Cells.Item(7, 47).Value = 7.5
$ws.Cells.Item(7, 48).Value = 41
$ws.Cells.Item(7, 49).Value = 6.5
$ws.Cells.Item(7, 50).Value = 21
$ws.Cells.Item(7, 51).Value = 26
$ws.Cells.Item(7, 52).Value = 67
$ws.Cells.Item(7, 53).Value = 81
$ws.Cells.Item(7, 54).Value = 126
$ws.Cells.Item(7, 55).Value = 351
$ws.Cells.Item(7, 56).Value = 81

# Row 8
$ws.Cells.Item(8, 1).Value = 'CpD7BeKb'
$ws.Cells.Item(8, 2).Value = '24/10/2024'
$ws.Cells.Item(8, 3).Value = '14:00'
$ws.Cells.Item(8, 4).Value = 'SPAIN - LALIGA2'
$ws.Cells.Item(8, 5).Value = 'Levante'
$ws.Cells.Item(8, 6).Value = 'Dep. La Coruna'
$ws.Cells.Item(8, 7).Value = 2.1
$ws.Cells.Item(8, 8).Value = 3.4
$ws.Cells.Item(8, 9).Value = 3.25
$ws.Cells.Item(8, 10).Value = 2.75
$ws.Cells.Item(8, 11).Value = 2.25
$ws.Cells.Item(8, 12).Value = 3.6
$ws.Cells.Item(8, 13).Value = 1.04
$ws.Cells.Item(8, 14).Value = 13
$ws.Cells.Item(8, 15).Value = 1.25
$ws.Cells.Item(8, 16).Value = 3.75
$ws.Cells.Item(8, 17).Value = 1.83
$ws.Cells.Item(8, 18).Value = 2.03
$ws.Cells.Item(8, 19).Value = 1.36
$ws.Cells.Item(8, 20).Value = 3
$ws.Cells.Item(8, 21).Value = 1.67
$ws.Cells.Item(8, 22).Value = 2.1
$ws.Cells.Item(8, 23).Value = 9
$ws.Cells.Item(8, 24).Value = 11
$ws.Cells.Item(8, 25).Value = 9
$ws.Cells.Item(8, 26).Value = 21
$ws.Cells.Item(8, 27).Value = 17
$ws.Cells.Item(8, 28).Value = 23
$ws.Cells.Item(8, 29).Value = 12
$ws.Cells.Item(8, 30).Value = 6.5
$ws.Cells.Item(8, 31).Value = 13
$ws.Cells.Item(8, 32).Value = 41
$ws.Cells.Item(8, 33).Value = 151
$ws.Cells.Item(8, 34).Value = 11
$ws.Cells.Item(8, 35).Value = 17
$ws.Cells.Item(8, 36).Value = 12
$ws.Cells.Item(8, 37).Value = 34
$ws.Cells.Item(8, 38).Value = 23
$ws.Cells.Item(8, 39).Value = 29
$ws.Cells.Item(8, 40).Value = 4.33
$ws.Cells.Item(8, 41).Value = 12
$ws.Cells.Item(8, 42).Value = 21
$ws.Cells.Item(8, 43).Value = 41
$ws.Cells.Item(8, 44).Value = 51
$ws.Cells.Item(8, 45).Value = 126
$ws.Cells.Item(8, 46).Value = 3
$ws.Cells.Item(8, 47).Value = 7.5
$ws.Cells.Item(8, 48).Value = 51
$ws.Cells.Item(8, 49).Value = 5
$ws.Cells.Item(8, 50).Value = 17
$ws.Cells.Item(8, 51).Value = 23
$ws.Cells.Item(8, 52).Value = 51
$ws.Cells.Item(8, 53).Value = 67
$ws.Cells.Item(8, 54).Value = 151
$ws.Cells.Item(8, 55).Value = 81
$ws.Cells.Item(8, 56).Value = 81
